# fix: removed duplicated lines
#
# The "Arquivo / Licença" table had two duplicated rows:
#   - row 22 ("Museu do Futebol" / CC-BY-SA) was an exact duplicate of row 15
#   - row 24 ("MUSEU DE ASTRONOMIA E CIÊNCIAS AFINS (II)" / CC-BY) was a
#     near-duplicate of row 23 ("MUSEU DE ASTRONOMIA E CIÊNCIAS AFINS" / CC-BY)
#
# Delete both duplicated rows outright (deleting from the bottom up so the
# remaining row indices don't shift while we work). Row 23 then slides up to
# become the new row 22, and the now-unused "(II)" string drops out of
# sharedStrings.xml automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(24).EntireRow.Delete()
$ws.Rows.Item(22).EntireRow.Delete()

# Match the saved view state: scrolled down so row 13 is at the top, with
# A22 (the last data row) selected.
$ws.Range("A13").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("A22").Select() | Out-Null
